$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("D16").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E16").Value = "['Normal']"

# Row 26
$ws.Range("D26").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E26").Value = "['HardwareFault']"

# Row 27
$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

# Row 28
$ws.Range("D28").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E28").Value = "['Normal', 'SoftwareFault']"

# Row 31
$ws.Range("D31").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E31").Value = "['Normal', 'SoftwareFault']"

# Row 35
$ws.Range("D35").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E35").Value = "['Normal', 'HardwareFault']"

# Row 36
$ws.Range("D36").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E36").Value = "['Normal', 'SurroundingEnvironment']"

# Row 38
$ws.Range("D38").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E38").Value = "['Normal', 'SoftwareFault']"

# Row 54
$ws.Range("D54").Value = "[0, 0, 0, 0, 0, 1, 0]"
$ws.Range("E54").Value = "['CommunicationIssue']"

# Row 61
$ws.Range("D61").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E61").Value = "['Normal', 'SoftwareFault']"

# Row 74
$ws.Range("D74").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E74").Value = "['Normal', 'SoftwareFault']"

# Row 81
$ws.Range("D81").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E81").Value = "['Normal']"

# Row 83
$ws.Range("D83").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E83").Value = "['Normal']"

# Row 84
$ws.Range("D84").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E84").Value = "['Normal']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal', 'HardwareFault']"

# Row 116
$ws.Range("D116").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E116").Value = "['Normal']"
